$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3851828718437105
$ws.Range("C2").Value = 0.04719667549206008
$ws.Range("D2").Value = 0.03340922047932793
$ws.Range("F2").Value = 0.7560649971391413
$ws.Range("G2").Value = 0.5978842922277039
$ws.Range("H2").Value = 0.7225724990014584
$ws.Range("K2").Value = 0.3558879256252396
$ws.Range("N2").Value = 1.394472939126404
$ws.Range("B3").Value = 0.3463348696886612
$ws.Range("C3").Value = 0.0420562393364321
$ws.Range("D3").Value = 0.03138271161257222
$ws.Range("F3").Value = 0.7560987703498228
$ws.Range("G3").Value = 0.5995907309549438
$ws.Range("H3").Value = 0.7276681375629295
$ws.Range("K3").Value = 0.3145961407293498
$ws.Range("N3").Value = 1.411309502454069
$ws.Range("B4").Value = 0.3225625386358502
$ws.Range("C4").Value = 0.0388811973444092
$ws.Range("D4").Value = 0.03012772173290301
$ws.Range("F4").Value = 0.7566338666132495
$ws.Range("G4").Value = 0.6010954284440473
$ws.Range("H4").Value = 0.7311544763514277
$ws.Range("K4").Value = 0.2892637444360844
$ws.Range("N4").Value = 1.422189524381263
$ws.Range("B5").Value = 0.3128957747833283
$ws.Range("C5").Value = 0.03758265852414411
$ws.Range("D5").Value = 0.029613640904536
$ws.Range("F5").Value = 0.7569811729353617
$ws.Range("G5").Value = 0.6018233213239768
$ws.Range("H5").Value = 0.7326651128339421
$ws.Range("K5").Value = 0.2789462549932864
$ws.Range("N5").Value = 1.42675954768505
$ws.Range("B6").Value = 0.3112918756095837
$ws.Range("C6").Value = 0.03736675552404733
$ws.Range("D6").Value = 0.02952811825511503
$ws.Range("F6").Value = 0.7570466466670567
$ws.Range("G6").Value = 0.6019511105811617
$ws.Range("H6").Value = 0.7329213842771694
$ws.Range("K6").Value = 0.2772333980130384
$ws.Range("N6").Value = 1.427526629181634
$ws.Range("B7").Value = 0.3224320851041966
$ws.Range("C7").Value = 0.03886370370342718
$ws.Range("D7").Value = 0.03012079939917811
$ws.Range("F7").Value = 0.7566380272978108
$ws.Range("G7").Value = 0.6011047808242296
$ws.Range("H7").Value = 0.7311744852144315
$ws.Range("K7").Value = 0.2891245756789544
$ws.Range("N7").Value = 1.422250605389106
$ws.Range("B8").Value = 0.3717716079278546
$ws.Range("C8").Value = 0.04542817758625972
$ws.Range("D8").Value = 0.03271271930164232
$ws.Range("F8").Value = 0.7559698305817548
$ws.Range("G8").Value = 0.5983777227826153
$ws.Range("H8").Value = 0.7242552584460498
$ws.Range("K8").Value = 0.341646405524358
$ws.Range("N8").Value = 1.400165607925306
$ws.Range("B9").Value = 0.4691518145274074
$ws.Range("C9").Value = 0.05815070039616899
$ws.Range("D9").Value = 0.03770944348285354
$ws.Range("F9").Value = 0.7587459197177182
$ws.Range("G9").Value = 0.5966646108615379
$ws.Range("H9").Value = 0.7135241037649109
$ws.Range("K9").Value = 0.4447948646547673
$ws.Range("N9").Value = 1.361161123346418
$ws.Range("B10").Value = 0.5410679625573209
$ws.Range("C10").Value = 0.06740540999074085
$ws.Range("D10").Value = 0.04132696727250362
$ws.Range("F10").Value = 0.7632856007728321
$ws.Range("G10").Value = 0.5976353558816356
$ws.Range("H10").Value = 0.7073701825558345
$ws.Range("K10").Value = 0.5206618939738803
$ws.Range("N10").Value = 1.335129157719072
$ws.Range("B11").Value = 0.5738632381134039
$ws.Range("C11").Value = 0.071595398280067
$ws.Range("D11").Value = 0.04296082115263289
$ws.Range("F11").Value = 0.765895706322965
$ws.Range("G11").Value = 0.598564046568967
$ws.Range("H11").Value = 0.7049464614474061
$ws.Range("K11").Value = 0.5551923891478339
$ws.Range("N11").Value = 1.323856623397912
$ws.Range("B12").Value = 0.5862931773230855
$ws.Range("C12").Value = 0.07317912695855
$ws.Range("D12").Value = 0.04357779980625764
$ws.Range("F12").Value = 0.7669625838736636
$ws.Range("G12").Value = 0.5989859830099107
$ws.Range("H12").Value = 0.7040827011807238
$ws.Range("K12").Value = 0.5682704967931045
$ws.Range("N12").Value = 1.319669960309147
$ws.Range("B13").Value = 0.5836156820700751
$ws.Range("C13").Value = 0.0728381737012711
$ws.Range("D13").Value = 0.04344499971818294
$ws.Range("F13").Value = 0.7667293200709224
$ws.Range("G13").Value = 0.5988919830905957
$ws.Range("H13").Value = 0.7042663231686817
$ws.Range("K13").Value = 0.5654538036269514
$ws.Range("N13").Value = 1.320567984439617
$ws.Range("B14").Value = 0.5748856372747184
$ws.Range("C14").Value = 0.07172575157501626
$ws.Range("D14").Value = 0.04301161514501928
$ws.Range("F14").Value = 0.7659819052991139
$ws.Range("G14").Value = 0.5985973501448001
$ws.Range("H14").Value = 0.7048743159999731
$ws.Range("K14").Value = 0.5562682923056457
$ws.Range("N14").Value = 1.32351054008809
$ws.Range("B15").Value = 0.5695396620257327
$ws.Range("C15").Value = 0.07104397800232221
$ws.Range("D15").Value = 0.04274592877990102
$ws.Range("F15").Value = 0.7655343170544597
$ws.Range("G15").Value = 0.5984260355240991
$ws.Range("H15").Value = 0.705253768740306
$ws.Range("K15").Value = 0.5506421688681371
$ws.Range("N15").Value = 1.325323622421079
$ws.Range("B16").Value = 0.5389262897916751
$ws.Range("C16").Value = 0.06713117722304673
$ws.Range("D16").Value = 0.04121995129652589
$ws.Range("F16").Value = 0.7631260009013303
$ws.Range("G16").Value = 0.5975844839838942
$ws.Range("H16").Value = 0.7075361408381582
$ws.Range("K16").Value = 0.5184055719682021
$ws.Range("N16").Value = 1.335877313571251
$ws.Range("B17").Value = 0.5201662059066336
$ws.Range("C17").Value = 0.06472563657227681
$ws.Range("D17").Value = 0.0402807739600064
$ws.Range("F17").Value = 0.7617882451964917
$ws.Range("G17").Value = 0.5971931311297567
$ws.Range("H17").Value = 0.709032550133486
$ws.Range("K17").Value = 0.4986337996756163
$ws.Range("N17").Value = 1.342497605303754
$ws.Range("B18").Value = 0.5093834838041573
$ws.Range("C18").Value = 0.06334015126404324
$ws.Range("D18").Value = 0.03973947756755081
$ws.Range("F18").Value = 0.761070092682381
$ws.Range("G18").Value = 0.5970138699255045
$ws.Range("H18").Value = 0.7099286092477115
$ws.Range("K18").Value = 0.4872633400198936
$ws.Range("N18").Value = 1.346359023267631
$ws.Range("B19").Value = 0.505733955621821
$ws.Range("C19").Value = 0.06287072802152238
$ws.Range("D19").Value = 0.03955601488681992
$ws.Range("F19").Value = 0.7608357438929971
$ws.Range("G19").Value = 0.5969610406477841
$ws.Range("H19").Value = 0.710238072646419
$ws.Range("K19").Value = 0.4834138107190711
$ws.Range("N19").Value = 1.347675636044025
$ws.Range("B20").Value = 0.5221624674957184
$ws.Range("C20").Value = 0.06498190561369199
$ws.Range("D20").Value = 0.04038086574193045
$ws.Range("F20").Value = 0.7619253424676486
$ws.Range("G20").Value = 0.5972300460175006
$ws.Range("H20").Value = 0.7088695946087853
$ws.Range("K20").Value = 0.5007383617446237
$ws.Range("N20").Value = 1.341787315683783
$ws.Range("B21").Value = 0.5774495660170089
$ws.Range("C21").Value = 0.07205257680190869
$ws.Range("D21").Value = 0.04313895781415766
$ws.Range("F21").Value = 0.7661993080667173
$ws.Range("G21").Value = 0.5986819823754104
$ws.Range("H21").Value = 0.7046942666885201
$ws.Range("K21").Value = 0.5589662447198691
$ws.Range("N21").Value = 1.322644013853003
$ws.Range("B22").Value = 0.6136472499030958
$ws.Range("C22").Value = 0.07665657412995586
$ws.Range("D22").Value = 0.04493145266142307
$ws.Range("F22").Value = 0.7694501300638166
$ws.Range("G22").Value = 0.6000405446515202
$ws.Range("H22").Value = 0.7022805051010721
$ws.Range("K22").Value = 0.5970339004962852
$ws.Range("N22").Value = 1.310610773533121
$ws.Range("B23").Value = 0.5943221274967527
$ws.Range("C23").Value = 0.07420091459050582
$ws.Range("D23").Value = 0.04397569742967988
$ws.Range("F23").Value = 0.7676732011481988
$ws.Range("G23").Value = 0.5992779007143554
$ws.Range("H23").Value = 0.7035399403679463
$ws.Range("K23").Value = 0.5767154962776715
$ws.Range("N23").Value = 1.316989367545462
$ws.Range("B24").Value = 0.521259949187538
$ws.Range("C24").Value = 0.06486605418915303
$ws.Range("D24").Value = 0.04033561845225364
$ws.Range("F24").Value = 0.7618632021222851
$ws.Range("G24").Value = 0.5972132143584759
$ws.Range("H24").Value = 0.7089431554315411
$ws.Range("K24").Value = 0.4997868998043771
$ws.Range("N24").Value = 1.342108265294405
$ws.Range("B25").Value = 0.4427419971639779
$ws.Range("C25").Value = 0.05472510881656945
$ws.Range("D25").Value = 0.03636702338749132
$ws.Range("F25").Value = 0.757556512471723
$ws.Range("G25").Value = 0.5967375323810842
$ws.Range("H25").Value = 0.7161233506380995
$ws.Range("K25").Value = 0.4168750526334577
$ws.Range("N25").Value = 1.37125202053236
